# Applies the "Added police data analysis draft" edit:
#  - Clears placeholder N/R and N/C cells that shouldn't have been populated
#  - Adds "Location" header label to two section header rows (42 and 50)
#  - Moves the active selection to F34

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear stray placeholder cells (N/R / N/C) ---------------------------

# Row 2
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# Row 3
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# Row 5
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()

# Row 6
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()

# Row 17
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()

# Row 18
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()

# Row 33
$ws.Range("C33").ClearContents()
$ws.Range("D33").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("F33").ClearContents()

# Row 34
$ws.Range("C34").ClearContents()
$ws.Range("D34").ClearContents()
$ws.Range("E34").ClearContents()
$ws.Range("F34").ClearContents()

# --- Add "Location" header labels -----------------------------------------

$ws.Range("B42").Value = "Location"
$ws.Range("B50").Value = "Location"

# --- Update active selection -----------------------------------------------

$ws.Range("F34").Select()
